# cv122042a.xlsx - "correção nos dados e inicio da analise PNAD 2009"
#
# The original sheet had two section-header rows ("situação do domicílio" at
# row 5 and "grandes regiões e unidades da federação" at row 8) that carried
# a label but no data, which had thrown off the alignment between the row
# labels (column A) and their data (columns B:I) for every row beneath them.
# The fix removes those two spacer rows entirely so the data shifts up and
# lines back up with the correct label, and renames the stray pandas
# artifact header "unnamed: 1_level_1" (B2) to "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the mis-named sub-header in B2.
$ws.Range("B2").Value = "total"

# Remove the "situação do domicílio" spacer row (row 5). Removing it shifts
# every row below up by one, so what had been row 6 ("urbana" data) becomes
# row 5, etc.
$ws.Rows(5).EntireRow.Delete()

# After the first deletion, the "grandes regiões e unidades da federação"
# spacer row (originally row 8) now sits at row 7. Remove it too, shifting
# everything below up by one more row.
$ws.Rows(7).EntireRow.Delete()
